$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New "Day 32 / Day 33 / Day 34" hyperlinked filenames in columns C/D/E for
# rows 33-35, mirroring the existing rows above them (e.g. row 32).
# ---------------------------------------------------------------------------

$links = @(
    @{ Cell = "C33"; Day = "Day 32"; File = "Count Complete Tree Nodes.java" },
    @{ Cell = "D33"; Day = "Day 32"; File = "Sum of left leaves.java" },
    @{ Cell = "E33"; Day = "Day 32"; File = "Merge Two Binary Trees.java" },
    @{ Cell = "C34"; Day = "Day 33"; File = "Binary Tree Paths.java" },
    @{ Cell = "D34"; Day = "Day 33"; File = "Lowest Common Ancestor of a Binary Tree.java" },
    @{ Cell = "E34"; Day = "Day 33"; File = "Root Equals Sum of Children.java" },
    @{ Cell = "C35"; Day = "Day 34"; File = "Univalued Binary Tree.java" }
)

foreach ($link in $links) {
    $target = "https://github.com/rhythmtaneja/90Day-LeetCode-Challenge/blob/main/" + $link.Day + "/" + $link.File
    $cell = $ws.Range($link.Cell)

    $ws.Hyperlinks.Add($cell, $target, "", $link.File, $target)
    $cell.Value = $link.File
    $cell.Style = "Hyperlink"
}

# ---------------------------------------------------------------------------
# Extend the trailing date-only rows (B39:B50) with three more days.
# Copy B50's formatting down so the new cells keep the same date style
# instead of minting a fresh one.
# ---------------------------------------------------------------------------

$newDates = @(
    @{ Cell = "B51"; Value = 45851 },
    @{ Cell = "B52"; Value = 45852 },
    @{ Cell = "B53"; Value = 45853 }
)

foreach ($d in $newDates) {
    $ws.Range("B50").Copy($ws.Range($d.Cell))
    $ws.Range($d.Cell).Value = $d.Value
}

# ---------------------------------------------------------------------------
# Match the author's final view state: scrolled a bit further down, with
# D40 selected.
# ---------------------------------------------------------------------------

$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D40").Select()
